$wb = $excel.ActiveWorkbook

# ALC!row17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1247.9375
$ws.Cells.Item(17, 10).Value = 1247.9375
$ws.Cells.Item(17, 12).Value = 3743.8125
$ws.Cells.Item(17, 14).Value = -4079.8125

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1499.9
$ws.Cells.Item(19, 9).Value = 2250
$ws.Cells.Item(19, 10).Value = 1312.375
$ws.Cells.Item(19, 11).Value = 2250
$ws.Cells.Item(19, 12).Value = 1312.375
$ws.Cells.Item(19, 13).Value = -2075
$ws.Cells.Item(19, 14).Value = -1662.375

# ALC!row51
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 9933.267
$ws.Cells.Item(51, 9).Value = 14975
$ws.Cells.Item(51, 10).Value = 4171.2856
$ws.Cells.Item(51, 11).Value = 14975
$ws.Cells.Item(51, 12).Value = 4171.2856
$ws.Cells.Item(51, 13).Value = -14491
$ws.Cells.Item(51, 14).Value = -5139.2856

# ALC!row53
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 454.48
$ws.Cells.Item(53, 9).Value = 541.25
$ws.Cells.Item(53, 10).Value = 413.64706
$ws.Cells.Item(53, 11).Value = 541.25
$ws.Cells.Item(53, 12).Value = 413.64706
$ws.Cells.Item(53, 13).Value = 95.75
$ws.Cells.Item(53, 14).Value = -1687.64706

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 341000
$ws.Cells.Item(86, 9).Value = 501500
$ws.Cells.Item(86, 10).Value = 20000
$ws.Cells.Item(86, 11).Value = 501500
$ws.Cells.Item(86, 12).Value = 20000
$ws.Cells.Item(86, 13).Value = -500377
$ws.Cells.Item(86, 14).Value = -22246

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 341000
$ws.Cells.Item(89, 9).Value = 501500
$ws.Cells.Item(89, 10).Value = 20000
$ws.Cells.Item(89, 11).Value = 2507500
$ws.Cells.Item(89, 12).Value = 100000
$ws.Cells.Item(89, 13).Value = -2501884
$ws.Cells.Item(89, 14).Value = -111232

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2984.75
$ws.Cells.Item(116, 9).Value = 2519.8
$ws.Cells.Item(116, 10).Value = 3759.6667
$ws.Cells.Item(116, 11).Value = 2519.8
$ws.Cells.Item(116, 12).Value = 3759.6667
$ws.Cells.Item(116, 13).Value = 922.1999999999998
$ws.Cells.Item(116, 14).Value = -10643.6667

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 723.6111
$ws.Cells.Item(74, 9).Value = 705.7857
$ws.Cells.Item(74, 10).Value = 786
$ws.Cells.Item(74, 11).Value = 705.7857
$ws.Cells.Item(74, 12).Value = 786
$ws.Cells.Item(74, 13).Value = 168.2143
$ws.Cells.Item(74, 14).Value = -2534

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 723.6111
$ws.Cells.Item(77, 9).Value = 705.7857
$ws.Cells.Item(77, 10).Value = 786
$ws.Cells.Item(77, 11).Value = 3528.9285
$ws.Cells.Item(77, 12).Value = 3930
$ws.Cells.Item(77, 13).Value = 839.0715
$ws.Cells.Item(77, 14).Value = -12666

# BSM!row50
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(50, 8).Value = 18730
$ws.Cells.Item(50, 10).Value = 18730
$ws.Cells.Item(50, 12).Value = 18730
$ws.Cells.Item(50, 14).Value = -19878

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3593.3428
$ws.Cells.Item(134, 9).Value = 3248.8215
$ws.Cells.Item(134, 10).Value = 4971.4287
$ws.Cells.Item(134, 11).Value = 9746.4645
$ws.Cells.Item(134, 12).Value = 14914.2861
$ws.Cells.Item(134, 13).Value = -7211.4645
$ws.Cells.Item(134, 14).Value = -19984.2861

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 26933.176
$ws.Cells.Item(31, 9).Value = 1524.2778
$ws.Cells.Item(31, 10).Value = 38660.36
$ws.Cells.Item(31, 11).Value = 1524.2778
$ws.Cells.Item(31, 12).Value = 38660.36
$ws.Cells.Item(31, 13).Value = -1229.2778
$ws.Cells.Item(31, 14).Value = -39250.36

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 26933.176
$ws.Cells.Item(34, 9).Value = 1524.2778
$ws.Cells.Item(34, 10).Value = 38660.36
$ws.Cells.Item(34, 11).Value = 1524.2778
$ws.Cells.Item(34, 12).Value = 38660.36
$ws.Cells.Item(34, 13).Value = -1322.2778
$ws.Cells.Item(34, 14).Value = -39064.36

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2381.6365
$ws.Cells.Item(62, 9).Value = 2116.5
$ws.Cells.Item(62, 10).Value = 2699.8
$ws.Cells.Item(62, 11).Value = 2116.5
$ws.Cells.Item(62, 12).Value = 2699.8
$ws.Cells.Item(62, 13).Value = -1492.5
$ws.Cells.Item(62, 14).Value = -3947.8  # new cell N62

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 2381.6365
$ws.Cells.Item(65, 9).Value = 2116.5
$ws.Cells.Item(65, 10).Value = 2699.8
$ws.Cells.Item(65, 11).Value = 10582.5
$ws.Cells.Item(65, 12).Value = 13499
$ws.Cells.Item(65, 13).Value = -7462.5
$ws.Cells.Item(65, 14).Value = -19739  # new cell N65

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 489.57144
$ws.Cells.Item(122, 9).Value = 345.5
$ws.Cells.Item(122, 10).Value = 681.6667
$ws.Cells.Item(122, 11).Value = 1036.5
$ws.Cells.Item(122, 12).Value = 2045.0001
$ws.Cells.Item(122, 13).Value = 1413.5
$ws.Cells.Item(122, 14).Value = -6945.0001

# CUL!row7
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 27
$ws.Cells.Item(7, 9).Value = 29.333334
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 88.00000199999999
$ws.Cells.Item(7, 12).Value = 60
$ws.Cells.Item(7, 13).Value = 23.99999800000001
$ws.Cells.Item(7, 14).Value = -284

# CUL!row33
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 2072.5
$ws.Cells.Item(33, 9).Value = 46.25
$ws.Cells.Item(33, 10).Value = 3423.3333
$ws.Cells.Item(33, 11).Value = 277.5
$ws.Cells.Item(33, 12).Value = 20539.9998
$ws.Cells.Item(33, 13).Value = 5.5
$ws.Cells.Item(33, 14).Value = -21105.9998

# CUL!row40
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 343.25
$ws.Cells.Item(40, 9).Value = 106.28571
$ws.Cells.Item(40, 11).Value = 425.14284
$ws.Cells.Item(40, 13).Value = -356.14284

# CUL!row44
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(44, 8).Value = 690.7
$ws.Cells.Item(44, 9).Value = 189
$ws.Cells.Item(44, 10).Value = 1025.1666
$ws.Cells.Item(44, 11).Value = 567
$ws.Cells.Item(44, 12).Value = 3075.4998
$ws.Cells.Item(44, 13).Value = -169
$ws.Cells.Item(44, 14).Value = -3871.4998

# CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1350
$ws.Cells.Item(68, 10).Value = 2300
$ws.Cells.Item(68, 12).Value = 6900
$ws.Cells.Item(68, 14).Value = -8522

# CUL!row69
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 1370
$ws.Cells.Item(69, 9).Value = 1980
$ws.Cells.Item(69, 10).Value = 1166.6666
$ws.Cells.Item(69, 11).Value = 5940
$ws.Cells.Item(69, 12).Value = 3499.9998
$ws.Cells.Item(69, 13).Value = -5129
$ws.Cells.Item(69, 14).Value = -5121.9998  # new cell N69

# CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 1350
$ws.Cells.Item(71, 10).Value = 2300
$ws.Cells.Item(71, 12).Value = 20700
$ws.Cells.Item(71, 14).Value = -28812

# CUL!row72
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 1370
$ws.Cells.Item(72, 9).Value = 1980
$ws.Cells.Item(72, 10).Value = 1166.6666
$ws.Cells.Item(72, 11).Value = 17820
$ws.Cells.Item(72, 12).Value = 10499.9994
$ws.Cells.Item(72, 13).Value = -13764
$ws.Cells.Item(72, 14).Value = -18611.9994  # new cell N72

# CUL!row80
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 13650.375
$ws.Cells.Item(80, 10).Value = 13650.375
$ws.Cells.Item(80, 12).Value = 40951.125
$ws.Cells.Item(80, 14).Value = -42823.125

# CUL!row83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 13650.375
$ws.Cells.Item(83, 10).Value = 13650.375
$ws.Cells.Item(83, 12).Value = 122853.375
$ws.Cells.Item(83, 14).Value = -132213.375

# CUL!row86
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 514.75
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 514.75
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 1544.25
$ws.Cells.Item(86, 14).Value = -3916.25
$ws.Cells.Item(86, 13).ClearContents()  # remove cell M86

# CUL!row89
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(89, 8).Value = 514.75
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 514.75
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 4632.75
$ws.Cells.Item(89, 14).Value = -16488.75
$ws.Cells.Item(89, 13).ClearContents()  # remove cell M89

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1348.8667
$ws.Cells.Item(113, 9).Value = 370.75
$ws.Cells.Item(113, 11).Value = 370.75
$ws.Cells.Item(113, 13).Value = 1799.25

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3357.2856
$ws.Cells.Item(126, 9).Value = 3164.8572
$ws.Cells.Item(126, 10).Value = 3742.1428
$ws.Cells.Item(126, 11).Value = 9494.571599999999
$ws.Cells.Item(126, 12).Value = 11226.4284
$ws.Cells.Item(126, 13).Value = -7024.571599999999
$ws.Cells.Item(126, 14).Value = -16166.4284

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1502.3334
$ws.Cells.Item(82, 9).Value = 500
$ws.Cells.Item(82, 10).Value = 1627.625
$ws.Cells.Item(82, 11).Value = 500
$ws.Cells.Item(82, 12).Value = 1627.625
$ws.Cells.Item(82, 13).Value = -139
$ws.Cells.Item(82, 14).Value = -2349.625

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1502.3334
$ws.Cells.Item(85, 9).Value = 500
$ws.Cells.Item(85, 10).Value = 1627.625
$ws.Cells.Item(85, 11).Value = 500
$ws.Cells.Item(85, 12).Value = 1627.625
$ws.Cells.Item(85, 13).Value = 748
$ws.Cells.Item(85, 14).Value = -4123.625

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2102.5
$ws.Cells.Item(136, 9).Value = 1642.3846
$ws.Cells.Item(136, 10).Value = 2957
$ws.Cells.Item(136, 11).Value = 4927.1538
$ws.Cells.Item(136, 12).Value = 8871
$ws.Cells.Item(136, 13).Value = -2377.1538
$ws.Cells.Item(136, 14).Value = -13971

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2154.8572
$ws.Cells.Item(122, 9).Value = 1680.6
$ws.Cells.Item(122, 10).Value = 2418.3333
$ws.Cells.Item(122, 11).Value = 5041.799999999999
$ws.Cells.Item(122, 12).Value = 7254.999899999999
$ws.Cells.Item(122, 13).Value = -2591.799999999999
$ws.Cells.Item(122, 14).Value = -12154.9999
